$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing "Employee Billing Utilization" entry ---
# Clear the previous cell formatting on the data cells first so the new
# number formats start from a clean (default) style instead of inheriting
# old font/alignment tweaks.
$ws.Range("A2:H2").ClearFormats()

$ws.Range("A2").Value = "sdasa"

$ws.Range("B2").Value = 42005
$ws.Range("B2").NumberFormat = "mm-dd-yy"

$ws.Range("C2").Value = 42369
$ws.Range("C2").NumberFormat = "mm-dd-yy"

$ws.Range("D2").Value = "ase"

$ws.Range("E2").Value = 0.0856
$ws.Range("E2").NumberFormat = "0.00%"

$ws.Range("F2").Value = 0.202
$ws.Range("F2").NumberFormat = "0.00%"

$ws.Range("G2").Value = 1
$ws.Range("G2").NumberFormat = "0.00%"

$ws.Range("H2").Value = 11
$ws.Range("H2").NumberFormat = "0"

# --- Row 3: new entry for Frank Warnakula ---
$ws.Range("A3").Value = "Frank Warnakula"

$ws.Range("B3").Value = 42006
$ws.Range("B3").NumberFormat = "mm-dd-yy"

$ws.Range("C3").Value = 42370
$ws.Range("C3").NumberFormat = "mm-dd-yy"

$ws.Range("D3").Value = "ase"

$ws.Range("E3").Value = 0.0856
$ws.Range("E3").NumberFormat = "0.00%"

$ws.Range("F3").Value = 0.202
$ws.Range("F3").NumberFormat = "0.00%"

$ws.Range("G3").Value = 1
$ws.Range("G3").NumberFormat = "0.00%"

$ws.Range("H3").Value = 1

# --- Row 4: new entry (no employee name) ---
$ws.Range("B4").Value = 42006
$ws.Range("B4").NumberFormat = "mm-dd-yy"

$ws.Range("C4").Value = 42370
$ws.Range("C4").NumberFormat = "mm-dd-yy"

$ws.Range("D4").Value = "ase"

$ws.Range("E4").Value = 0.0526
$ws.Range("E4").NumberFormat = "0.00%"

$ws.Range("F4").Value = 0.125
$ws.Range("F4").NumberFormat = "0.00%"

$ws.Range("G4").Value = 1
$ws.Range("G4").NumberFormat = "0.00%"

$ws.Range("H4").Value = 1

# --- Final selection, matching the saved workbook state ---
$ws.Range("H6").Select()
